# Incident View Form Test
# Adds a new "isReallyIncident" variable row to the Variables sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "isReallyIncident"
$ws.Range("B16").Value = "Is this really an Incident?"
$ws.Range("C16").Value = "Boolean"
$ws.Range("C16").WrapText = $true
$ws.Range("D16").Value = "Company"

# Leave selection on the next empty row, as in the authored workbook.
$ws.Range("A17").Select()
